$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), using the same formatting
# as the other header cells in row 1 (bold, centered, thin border).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: column I (I0) is always 1, column J (IF) mirrors column H's
# (IP) value for that row.
for ($row = 2; $row -le 32; $row++) {
    $hVal = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hVal
}
